# Update the "Fitness" values (column C) for rows 2-60 on the active sheet
# (Sheet1) to reflect the new run results, as captured by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C3").Value = 12184
$ws.Range("C4:C6").Value = 10706
$ws.Range("C7:C12").Value = 10427
$ws.Range("C13:C14").Value = 10390
$ws.Range("C15").Value = 9958
$ws.Range("C16:C17").Value = 9162
$ws.Range("C18:C22").Value = 8759
$ws.Range("C23:C25").Value = 8667
$ws.Range("C26:C30").Value = 8664
$ws.Range("C31:C37").Value = 8246
$ws.Range("C38:C43").Value = 7345
$ws.Range("C44:C53").Value = 7343
$ws.Range("C54:C60").Value = 7310
